$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "300.78"
Set-TextValue $ws.Range("E2") "-0.04%"
Set-TextValue $ws.Range("D3") "31.73"
Set-TextValue $ws.Range("E3") "1.07%"
Set-TextValue $ws.Range("D4") "5.090"
Set-TextValue $ws.Range("E4") "-1.09%"
Set-TextValue $ws.Range("D5") "0.08204"
Set-TextValue $ws.Range("E5") "11.07%"
Set-TextValue $ws.Range("D6") "2.622"
Set-TextValue $ws.Range("E6") "5.59%"
Set-TextValue $ws.Range("D7") "7.830"
Set-TextValue $ws.Range("E7") "-1.34%"
Set-TextValue $ws.Range("D9") "0.9271"
Set-TextValue $ws.Range("E9") "0.71%"
Set-TextValue $ws.Range("D10") "0.1752"
Set-TextValue $ws.Range("E10") "1.01%"
Set-TextValue $ws.Range("D11") "0.07461"
Set-TextValue $ws.Range("E11") "-0.60%"
Set-TextValue $ws.Range("D12") "0.08909"
Set-TextValue $ws.Range("E12") "9.55%"
Set-TextValue $ws.Range("D13") "0.02992"
Set-TextValue $ws.Range("E13") "-1.43%"
Set-TextValue $ws.Range("D14") "0.1004"
Set-TextValue $ws.Range("E14") "1.21%"
Set-TextValue $ws.Range("D15") "0.001522"
Set-TextValue $ws.Range("E15") "1.87%"
Set-TextValue $ws.Range("D16") "0.005745"
Set-TextValue $ws.Range("D17") "3.591"
Set-TextValue $ws.Range("E17") "3.78%"
Set-TextValue $ws.Range("D18") "2.259"
Set-TextValue $ws.Range("E18") "1.34%"
Set-TextValue $ws.Range("D19") "0.3245"
Set-TextValue $ws.Range("E19") "-1.02%"
Set-TextValue $ws.Range("D20") "0.1345"
Set-TextValue $ws.Range("E20") "0.50%"
Set-TextValue $ws.Range("D21") "3.894"
Set-TextValue $ws.Range("E21") "-16.21%"
Set-TextValue $ws.Range("D22") "0.1694"
Set-TextValue $ws.Range("E22") "7.99%"
Set-TextValue $ws.Range("D23") "0.04606"
Set-TextValue $ws.Range("E23") "-1.01%"
Set-TextValue $ws.Range("D24") "0.001240"
Set-TextValue $ws.Range("E24") "1.24%"
Set-TextValue $ws.Range("D25") "0.004521"
Set-TextValue $ws.Range("E25") "0.63%"
Set-TextValue $ws.Range("D26") "0.0001194"
Set-TextValue $ws.Range("E26") "-8.09%"
Set-TextValue $ws.Range("D27") "0.0003398"
Set-TextValue $ws.Range("E27") "81.57%"
Set-TextValue $ws.Range("D39") "0.01779"
Set-TextValue $ws.Range("E39") "2.90%"
Set-TextValue $ws.Range("D40") "0.04547"
Set-TextValue $ws.Range("E40") "0.65%"
Set-TextValue $ws.Range("D41") "0.006964"
Set-TextValue $ws.Range("E41") "-3.00%"
Set-TextValue $ws.Range("D42") "0.1375"
Set-TextValue $ws.Range("E42") "2.21%"
Set-TextValue $ws.Range("D43") "0.002132"
Set-TextValue $ws.Range("E43") "-4.34%"
Set-TextValue $ws.Range("D44") "0.009584"
Set-TextValue $ws.Range("E44") "-10.86%"
Set-TextValue $ws.Range("D45") "0.00006461"
Set-TextValue $ws.Range("E45") "2.97%"
Set-TextValue $ws.Range("E46") "-0.48%"
Set-TextValue $ws.Range("E47") "-12.63%"
Set-TextValue $ws.Range("D48") "0.8206"
Set-TextValue $ws.Range("E48") "-57.45%"
Set-TextValue $ws.Range("D49") "0.00002090"
Set-TextValue $ws.Range("E49") "-0.48%"
Set-TextValue $ws.Range("D50") "0.0001990"
Set-TextValue $ws.Range("E50") "-0.41%"
